$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 60848.79486621871
$ws.Range("D2").Value = 91374816049.77301
$ws.Range("G2").Value = 85015604973.8835

$ws.Range("B3").Value = 61436.25340635831
$ws.Range("D3").Value = 91348500629.70424
$ws.Range("G3").Value = 83455906330.53168

$ws.Range("B4").Value = 68704.53917792696
$ws.Range("D4").Value = 91344620460.91911
$ws.Range("G4").Value = 84964515780.99345

$ws.Range("B5").Value = 69288.90071436594
$ws.Range("D5").Value = 91315830551.4501
$ws.Range("G5").Value = 83318508790.02858

$ws.Range("B6").Value = 56337.08762636479
$ws.Range("D6").Value = 84554801714.51651
$ws.Range("G6").Value = 15527131851.68781

$ws.Range("B7").Value = 56874.21707019051
$ws.Range("D7").Value = 84526998498.32367
$ws.Range("G7").Value = 13776448083.28479

$ws.Range("B8").Value = 56024.38565162199
$ws.Range("D8").Value = 84449457716.90875
$ws.Range("G8").Value = 15758505557.68437

$ws.Range("B9").Value = 56564.89569795261
$ws.Range("D9").Value = 84421950578.89888
$ws.Range("G9").Value = 13995665498.02525
